# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" worksheet (cloned from the "2021-Q4" layout) right
# before the "总计" sheet, fills it in with the Q1-2022 fund holding figures,
# and prepends a matching summary row to the "总计" sheet.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q1" sheet, positioned right after "2021-Q4"
#    (i.e. right before "总计").
# ---------------------------------------------------------------------
$q4Sheet = $wb.Worksheets.Item("2021-Q4")
$q1Sheet = $wb.Worksheets.Add($null, $q4Sheet)
$q1Sheet.Name = "2022-Q1"

# Clone the header row + layout/styles from "2021-Q4" (same column headers:
# 基金代码/基金名称/基金规模/股票总仓位/仓位占比/持有市值(亿元)/仓位排名).
$q4Sheet.Range("A1:H2").Copy($q1Sheet.Range("A1"))
$q1Sheet.Range("A1").Clear()

# Overwrite the data row with the 2022-Q1 figures. Columns D:G hold
# text-looking numbers (kept as text, matching the source data), so force
# a text format before assigning, then drop back to the normal style so no
# stray number formatting is left behind.
$q1Sheet.Range("D2:G2").NumberFormat = "@"
$q1Sheet.Range("D2").Value = "6.05"
$q1Sheet.Range("E2").Value = "99.49"
$q1Sheet.Range("F2").Value = "4.19"
$q1Sheet.Range("G2").Value = "0.2535"
$q1Sheet.Range("D2:G2").Style = "Normal"
$q1Sheet.Range("H2").Value = 10

# ---------------------------------------------------------------------
# 2. Prepend a "2022-Q1" row to the "总计" summary sheet, pushing the
#    existing rows down.
# ---------------------------------------------------------------------
$totalSheet = $wb.Worksheets.Item("总计")
$totalSheet.Rows.Item(2).Insert()

# The row-insert copies the bold/bordered header formatting down into the
# new row; reset B2:D2 back to the plain style used by the other data rows.
$totalSheet.Range("B2:D2").Style = "Normal"

# A2 keeps the bordered "index" style used by the column-A cells; copy it
# from the row just below (which held that same style pre-insert).
$totalSheet.Range("A3").Copy($totalSheet.Range("A2"))

$totalSheet.Range("A2").Value = 0
$totalSheet.Range("B2").Value = "2022-Q1"
$totalSheet.Range("C2").Value = 1
$totalSheet.Range("D2").Value = 0.25

# Renumber the index column for the rows that shifted down.
$totalSheet.Range("A3").Value = 1
$totalSheet.Range("A4").Value = 2
